$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 becomes the end of its block: reuse the existing "end of block" formats
# (style indices 6/7 already defined in the workbook) by copying formats from
# row 3, which already uses that style. PasteSpecial(formats only) reuses the
# existing style entries instead of creating new ones.
$ws.Range("A3:B3").Copy()
$ws.Range("A14:B14").PasteSpecial(-4122)
$ws.Range("C3:E3").Copy()
$ws.Range("C14:E14").PasteSpecial(-4122)

# New row 15 gets the regular (non-end-of-block) formatting, copied from row 13
# which already uses style indices 4/5.
$ws.Range("A13:B13").Copy()
$ws.Range("A15:B15").PasteSpecial(-4122)
$ws.Range("C13:E13").Copy()
$ws.Range("C15:E15").PasteSpecial(-4122)

# Fill in the new row's content (order chosen to reproduce the original
# shared-string table ordering).
$ws.Range("C15").Value = " The guild seems to be buzzing\nwith activity…"
$ws.Range("A15").Value = "SCRIPT/G01P03A/um2102.ssb"
$ws.Range("B15").Value = 197
$ws.Range("D15").Value = " Члены гильдии сегодня очень\nактивны..."
$ws.Range("E15").Value = " Œìåîú ãéìûäéé òåãïäîÿ ïœåîû\nàëóéâîú..."

# Match the row height used by the other wrapped-text rows.
$ws.Rows.Item(15).RowHeight = 43.2

# Match the new selection recorded in the sheet view.
$ws.Range("C14").Select()
